# Fix contact information missing from short resumes:
# insert a centered contact-info line right after the "Dheeraj Chand" header
# paragraph, using a paragraph-mark (^p) replacement so the new run does not
# inherit the header run's Bold/Size formatting.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
